$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.809.32'
$ws.Range("E2").Value = '  +0.14%  '
# Row 3
$ws.Range("D3").Value = '2.290.49'
$ws.Range("E3").Value = '  -0.06%  '
# Row 4
$ws.Range("E4").Value = '  +0.13%  '
# Row 5
$ws.Range("D5").Value = '''115.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +17.23%  '
# Row 6
$ws.Range("D6").Value = '''269.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '
# Row 7
$ws.Range("D7").Value = '''0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.40%  '
# Row 8
$ws.Range("E8").Value = '  +0.24%  '
# Row 9
$ws.Range("D9").Value = '''0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.72%  '
# Row 10
$ws.Range("D10").Value = '''48.85'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.30%  '
# Row 11
$ws.Range("D11").Value = '''0.0944'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.15%  '
# Row 12
$ws.Range("D12").Value = '''9.05'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +14.26%  '
# Row 13
$ws.Range("D13").Value = '''0.108'
$ws.Range("D13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = '''15.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.83%  '
# Row 15
$ws.Range("D15").Value = '2.635.52'
$ws.Range("E15").Value = '  +0.01%  '
# Row 16
$ws.Range("D16").Value = '''0.879'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.10%  '
# Row 17
$ws.Range("D17").Value = '2.297.28'
$ws.Range("E17").Value = '  +0.15%  '
# Row 18
$ws.Range("D18").Value = '43.687.62'
$ws.Range("E18").Value = '  -0.10%  '
# Row 19
$ws.Range("D19").Value = '''0.0000110'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.18%  '
# Row 20
$ws.Range("D20").Value = '''7.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +12.78%  '
# Row 21
$ws.Range("D21").Value = '''72.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.47%  '
# Row 22
$ws.Range("E22").Value = '  -1.72%  '
# Row 23
$ws.Range("D23").Value = '''10.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.64%  '
# Row 24
$ws.Range("D24").Value = '''233.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.19%  '
# Row 25
$ws.Range("D25").Value = '''2.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.50%  '
# Row 26
$ws.Range("D26").Value = '''11.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.79%  '
# Row 27
$ws.Range("E27").Value = '  +0.00%  '
# Row 28
$ws.Range("E28").Value = '  +4.75%  '
# Row 29
$ws.Range("D29").Value = '''41.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.45%  '
# Row 30
$ws.Range("D30").Value = '''3.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.96%  '
# Row 31
$ws.Range("D31").Value = '''2.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.20%  '
# Row 32
$ws.Range("D32").Value = '''173.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.63%  '
# Row 33
$ws.Range("D33").Value = '''0.0937'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.16%  '
# Row 34
$ws.Range("E34").Value = '  -1.62%  '
# Row 35
$ws.Range("D35").Value = '''5.74'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.79%  '
# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '''0.128'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.13%  '
# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '''4.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.11%  '
# Row 38
$ws.Range("D38").Value = '''0.0362'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.12%  '
# Row 39
$ws.Range("D39").Value = '''0.108'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.24%  '
# Row 40
$ws.Range("D40").Value = '''3.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.17%  '
# Row 41
$ws.Range("D41").Value = '''14.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +18.42%  '
# Row 42
$ws.Range("D42").Value = '''75.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.42%  '
# Row 43
$ws.Range("D43").Value = '''2.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.64%  '
# Row 44
$ws.Range("D44").Value = '''0.243'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.65%  '
# Row 45
$ws.Range("D45").Value = '''6.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +22.49%  '
# Row 46
$ws.Range("E46").Value = '  +0.21%  '
# Row 47
$ws.Range("D47").Value = '''1.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.97%  '
# Row 48
$ws.Range("E48").Value = '  -0.95%  '
# Row 49
$ws.Range("D49").Value = '''103.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.42%  '
# Row 50
$ws.Range("E50").Value = '  +3.57%  '
# Row 51
$ws.Range("D51").Value = '''0.0999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
